# "Generate Report for Handback"
#
# The localization-status report is regenerated after a handback event:
#  - Overview + per-language "Status" cells flip from "Ready for handoff" to
#    "Handed back: in sync with en-US".
#  - Each language sheet (zh-cn, de-de) gets its "Latest Target File",
#    "Latest Handback File" and "Latest Handback DateTime" columns filled in
#    for both tracked source documents.
#  - A couple of columns are widened to fit the newly-populated values.

$wb = $excel.ActiveWorkbook

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99656b818048de01e5211edcb331f204ec2e7292/e2e/2c59594f-d292-4b4f-ada1-1b7351014f77.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99656b818048de01e5211edcb331f204ec2e7292/e2e/fd392e74-0161-411a-b6bf-256918985800.md"
$mdName1 = "2c59594f-d292-4b4f-ada1-1b7351014f77.md"
$mdName2 = "fd392e74-0161-411a-b6bf-256918985800.md"

$newStatus = "Handed back: in sync with en-US"

# ColumnWidth inputs chosen so the saved column width lands on the closest
# achievable value to the authored widths (29.9777047293527 -> 30,
# 18.6506053379604 / 21.7054770333426 -> 40).
$wideColWidth = 29.166666666666664
$fullColWidth = 39.16666666666667

# --- Overview sheet -------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus

$ov.Columns.Item(5).ColumnWidth = $wideColWidth
$ov.Columns.Item(6).ColumnWidth = $wideColWidth

# --- zh-cn sheet ------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl1, "", "", $mdName1)
$zh.Range("I2").Style = "HyperLink"
$zh.Range("J2").Value = "2c59594f-d292-4b4f-ada1-1b7351014f77.659f73a3078e4b7536f8e4a4c469556d96689450.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-17 22:59:49"

$zh.Hyperlinks.Add($zh.Range("I3"), $mdUrl2, "", "", $mdName2)
$zh.Range("I3").Style = "HyperLink"
$zh.Range("J3").Value = "fd392e74-0161-411a-b6bf-256918985800.3c3b001d6d4bc51ba5059cb70d58664fb7be567c.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-17 22:59:49"

$zh.Columns.Item(3).ColumnWidth = $wideColWidth
$zh.Columns.Item(9).ColumnWidth = $fullColWidth
$zh.Columns.Item(10).ColumnWidth = $fullColWidth

# --- de-de sheet ------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

$de.Hyperlinks.Add($de.Range("I2"), $mdUrl1, "", "", $mdName1)
$de.Range("I2").Style = "HyperLink"
$de.Range("J2").Value = "2c59594f-d292-4b4f-ada1-1b7351014f77.659f73a3078e4b7536f8e4a4c469556d96689450.de-de.xlf"
$de.Range("K2").Value = "2016-08-17 22:59:56"

$de.Hyperlinks.Add($de.Range("I3"), $mdUrl2, "", "", $mdName2)
$de.Range("I3").Style = "HyperLink"
$de.Range("J3").Value = "fd392e74-0161-411a-b6bf-256918985800.3c3b001d6d4bc51ba5059cb70d58664fb7be567c.de-de.xlf"
$de.Range("K3").Value = "2016-08-17 22:59:56"

$de.Columns.Item(3).ColumnWidth = $wideColWidth
$de.Columns.Item(9).ColumnWidth = $fullColWidth
$de.Columns.Item(10).ColumnWidth = $fullColWidth
